$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Generator Data" ---
$ws1 = $wb.Worksheets.Item("Generator Data")

# Relabel rows that used to describe "upgrade 2 / upgrade 3" tiers; the
# model now only has a single upgrade level, so the former "upgrade 2"
# text/values slide up into the "upgrade 1" follow-on rows.
$ws1.Range("A7").Value = "Investment at upgrade 1"
$ws1.Range("A8").Value = "Yearly O&M Cost at upgrade 1"
$ws1.Range("A9").Value = "Total actualized Fuel Cost"

# Updated numeric results for the new run.
$ws1.Range("B6").Value = 144214.78917984059
$ws1.Range("B7").Value = 60584.632934451052
$ws1.Range("B8").Value = 6058.4632934451056
$ws1.Range("B9").Value = 1130449.8563873509

# Remove the now-unused "upgrade 2" / "upgrade 3" rows entirely.
$ws1.Range("A10:A15").EntireRow.Delete()

# --- Sheet 2: "Yearly Fuel Costs" ---
$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")

# Updated numeric results for years 1-3.
$ws2.Range("B2").Value = 174013.2147760638
$ws2.Range("B3").Value = 198224.2912546304
$ws2.Range("B4").Value = 223257.02587713659

# The run now covers 5 years instead of 3, so append the two new rows,
# copying the label column's formatting (style index) from the row above.
$ws2.Range("A5").Value = "Total Fuel Cost at y = 4"
$ws2.Range("B5").Value = 252371.9633319019
$ws2.Range("A6").Value = "Total Fuel Cost at y = 5"
$ws2.Range("B6").Value = 282583.36114769487

$ws2.Range("A4").Copy()
$ws2.Range("A5:A6").PasteSpecial(-4122)
